$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 / TestScenario_4 precondition column picked up a stray single-space
# entry (matches the "Precondition" column, D) before the Reject edits below.
$ws.Range("D13").Value = " "

# Mark the two TestCase_1 rows (row 9 - TestScenario_3, row 14 - TestScenario_4)
# as Rejected in the Approved/Rejected column (I).
$ws.Range("I9").Value = "Rejected"
$ws.Range("I14").Value = "Rejected"

# Apply an AutoFilter on the "Steps" column (F, Field 6 within A1:K16) so only
# rows with Steps = "Step 1" remain visible (xlFilterValues = 7).
$ws.Range("A1:K16").AutoFilter(6, "Step 1", 7)

# Scroll the view over to show column H onward, with I14 as the active cell.
$ws.Range("I14").Select()
$ws.Application.ActiveWindow.ScrollColumn = 8
